# This script updates the cryptocurrency price/volume table with fresh
# values. Columns D (Price) and E (Volume(1h)) in the source data are
# plain text (not numbers), so for any D-column value that *looks* like
# a pure number (e.g. "0.999", "53.49") we force the cell to Text format
# before assigning, then restore the default "Normal" style so no stray
# number-format style sticks to the cell. Values that already contain
# multiple "." separators (e.g. "73.105.34") or other non-numeric
# characters never auto-convert, so they can be set directly.

function Set-TextValue($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
Set-TextValue $ws "D2" "73.105.34"
$ws.Range("E2").Value = "  +3.02%  "

# Row 3 - Ethereum
Set-TextValue $ws "D3" "3.984.50"
$ws.Range("E3").Value = "  +0.88%  "

# Row 4 - TetherUSD
Set-TextValue $ws "D4" "0.999"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5 - BNB
Set-TextValue $ws "D5" "597.90"
$ws.Range("E5").Value = "  +11.45%  "

# Row 6 - Solana
Set-TextValue $ws "D6" "159.89"
$ws.Range("E6").Value = "  +7.89%  "

# Row 7 - XRP
Set-TextValue $ws "D7" "0.682"
$ws.Range("E7").Value = "  -0.49%  "

# Row 9 - Cardano
Set-TextValue $ws "D9" "0.749"
$ws.Range("E9").Value = "  +1.80%  "

# Row 10 - Dogecoin (only E changes)
$ws.Range("E10").Value = "  +1.89%  "

# Row 11 - Avalanche
Set-TextValue $ws "D11" "53.49"
$ws.Range("E11").Value = "  -3.74%  "

# Row 12 - ShibaInu (only E changes)
$ws.Range("E12").Value = "  +0.88%  "

# Row 13 - Polkadot
Set-TextValue $ws "D13" "10.97"
$ws.Range("E13").Value = "  +3.47%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue $ws "D14" "4.613.50"
$ws.Range("E14").Value = "  +0.75%  "

# Row 15 - WrappedEther
Set-TextValue $ws "D15" "3.984.41"
$ws.Range("E15").Value = "  +0.58%  "

# Row 16 - Polygon (only E changes)
$ws.Range("E16").Value = "  +8.20%  "

# Row 17 - Uniswap
Set-TextValue $ws "D17" "14.06"
$ws.Range("E17").Value = "  +1.74%  "

# Row 18 - Chainlink
Set-TextValue $ws "D18" "20.31"
$ws.Range("E18").Value = "  -1.18%  "

# Row 19 - TRON (only E changes)
$ws.Range("E19").Value = "  +0.30%  "

# Row 20 - WrappedBTC
Set-TextValue $ws "D20" "72.675.78"
$ws.Range("E20").Value = "  +2.57%  "

# Row 21 - BitcoinCash
Set-TextValue $ws "D21" "434.51"
$ws.Range("E21").Value = "  +2.25%  "

# Row 22 - PancakeSwap
Set-TextValue $ws "D22" "4.79"
$ws.Range("E22").Value = "  +13.96%  "

# Row 23 - Litecoin
Set-TextValue $ws "D23" "95.93"
$ws.Range("E23").Value = "  -1.03%  "

# Row 24 - ImmutableX
Set-TextValue $ws "D24" "3.42"
$ws.Range("E24").Value = "  -4.36%  "

# Row 25 - InternetComputer(DFINITY)
Set-TextValue $ws "D25" "14.22"
$ws.Range("E25").Value = "  -0.94%  "

# Row 26 - Toncoin
Set-TextValue $ws "D26" "4.35"
$ws.Range("E26").Value = "  +15.58%  "

# Row 27 - RenderToken
Set-TextValue $ws "D27" "11.15"
$ws.Range("E27").Value = "  -2.17%  "

# Row 28 - LEO (only E changes)
$ws.Range("E28").Value = "  +1.05%  "

# Row 29 - Filecoin
Set-TextValue $ws "D29" "10.45"
$ws.Range("E29").Value = "  -1.46%  "

# Row 30 - EthereumClassic
Set-TextValue $ws "D30" "36.28"
$ws.Range("E30").Value = "  -0.07%  "

# Row 31 - NEARProtocol
Set-TextValue $ws "D31" "7.79"
$ws.Range("E31").Value = "  -0.46%  "

# Row 32 - Cosmos (only E changes)
$ws.Range("E32").Value = "  +2.90%  "

# Row 33 - Hedera (only E changes)
$ws.Range("E33").Value = "  +0.28%  "

# Row 34 - InjectiveProtocol
Set-TextValue $ws "D34" "47.81"
$ws.Range("E34").Value = "  -6.03%  "

# Row 35 - Bittensor
Set-TextValue $ws "D35" "665.16"
$ws.Range("E35").Value = "  -2.61%  "

# Row 36 - OKB
Set-TextValue $ws "D36" "70.89"
$ws.Range("E36").Value = "  +9.08%  "

# Row 37 - PEPE
Set-TextValue $ws "D37" "0.0₃0902"
$ws.Range("E37").Value = "  +10.10%  "

# Row 38 - TheGraph
Set-TextValue $ws "D38" "0.437"
$ws.Range("E38").Value = "  +0.22%  "

# Row 39 - Dai (only E changes)
$ws.Range("E39").Value = "  -0.06%  "

# Row 40 - now WEMIXToken (was ThetaToken)
$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws "D40" "3.35"
$ws.Range("E40").Value = "  +5.30%  "

# Row 41 - now ThetaToken (was WEMIXToken)
$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws "D41" "3.34"
$ws.Range("E41").Value = "  -1.16%  "

# Row 42 - Kaspa (only E changes)
$ws.Range("E42").Value = "  -2.88%  "

# Row 43 - FirstDigitalUSD (only E changes)
$ws.Range("E43").Value = "  +0.18%  "

# Row 44 - VeChain
Set-TextValue $ws "D44" "0.0489"
$ws.Range("E44").Value = "  +1.92%  "

# Row 45 - THORChain (only E changes)
$ws.Range("E45").Value = "  +6.62%  "

# Row 46 - Stellar (only E changes)
$ws.Range("E46").Value = "  +0.52%  "

# Row 47 - ApeXProtocol (only E changes)
$ws.Range("E47").Value = "  +3.33%  "

# Row 48 - Fetch.AI (only E changes)
$ws.Range("E48").Value = "  -3.39%  "

# Row 49 - Maker
Set-TextValue $ws "D49" "2.880.80"
$ws.Range("E49").Value = "  +9.81%  "

# Row 50 - Stacks
Set-TextValue $ws "D50" "3.03"
$ws.Range("E50").Value = "  +1.07%  "

# Row 51 - LidoDAOToken (only E changes)
$ws.Range("E51").Value = "  +4.37%  "
